$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.331.91"
$ws.Cells.Item(2, 5).Value = "  -0.08%  "
$ws.Cells.Item(3, 4).Value = "1.840.32"
$ws.Cells.Item(3, 5).Value = "  -0.26%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "239.08"
$ws.Cells.Item(5, 5).Value = "  -0.52%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.6289"
$ws.Cells.Item(6, 5).Value = "  +0.28%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.9998"
$ws.Cells.Item(7, 5).Value = "  +0.10%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.07432"
$ws.Cells.Item(8, 5).Value = "  -0.69%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "25.01"
$ws.Cells.Item(9, 5).Value = "  +2.66%  "
$ws.Cells.Item(10, 5).Value = "  -0.36%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07724"
$ws.Cells.Item(11, 5).Value = "  +0.14%  "
$ws.Cells.Item(12, 4).Value = "1.833.46"
$ws.Cells.Item(12, 5).Value = "  -0.63%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.951"
$ws.Cells.Item(13, 5).Value = "  -1.05%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.6738"
$ws.Cells.Item(14, 5).Value = "  -0.78%  "
$ws.Cells.Item(15, 5).Value = "  -0.74%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "81.64"
$ws.Cells.Item(16, 5).Value = "  -0.64%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "6.211"
$ws.Cells.Item(17, 5).Value = "  +0.73%  "
$ws.Cells.Item(18, 4).Value = "29.283.20"
$ws.Cells.Item(18, 5).Value = "  -0.34%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "229.10"
$ws.Cells.Item(19, 5).Value = "  +0.13%  "
$ws.Cells.Item(20, 5).Value = "  -0.51%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.9999"
$ws.Cells.Item(21, 5).Value = "  +0.14%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "7.337"
$ws.Cells.Item(22, 5).Value = "  -1.72%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.9999"
$ws.Cells.Item(23, 5).Value = "  +0.09%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "158.21"
$ws.Cells.Item(24, 5).Value = "  -0.47%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "8.466"
$ws.Cells.Item(25, 5).Value = "  +0.69%  "
$ws.Cells.Item(26, 5).Value = "  -2.35%  "
$ws.Cells.Item(27, 2).Value = "Hedera"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.07501"
$ws.Cells.Item(27, 5).Value = "  +17.50%  "
$ws.Cells.Item(28, 2).Value = "EthereumClassic"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "17.34"
$ws.Cells.Item(28, 5).Value = "  -1.16%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.455"
$ws.Cells.Item(29, 5).Value = "  +5.12%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.475"
$ws.Cells.Item(30, 5).Value = "  +0.41%  "
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.037"
$ws.Cells.Item(31, 5).Value = "  -1.41%  "
$ws.Cells.Item(32, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.038"
$ws.Cells.Item(32, 5).Value = "  -0.49%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.815"
$ws.Cells.Item(33, 5).Value = "  -0.54%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.139"
$ws.Cells.Item(34, 5).Value = "  -0.17%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.6903"
$ws.Cells.Item(35, 5).Value = "  -1.38%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.572"
$ws.Cells.Item(36, 5).Value = "  -0.17%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01840"
$ws.Cells.Item(37, 5).Value = "  +0.67%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "6.896"
$ws.Cells.Item(38, 5).Value = "  +4.64%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.811"
$ws.Cells.Item(39, 5).Value = "  -0.55%  "
$ws.Cells.Item(40, 4).Value = "1.233.74"
$ws.Cells.Item(40, 5).Value = "  -2.03%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.9328"
$ws.Cells.Item(41, 5).Value = "  +2.62%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.9992"
$ws.Cells.Item(42, 5).Value = "  +0.06%  "
$ws.Cells.Item(43, 2).Value = "Quant"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "101.05"
$ws.Cells.Item(43, 5).Value = "  -0.35%  "
$ws.Cells.Item(44, 2).Value = "RocketPoolETH"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(44, 4).Value = "1.973.25"
$ws.Cells.Item(44, 5).Value = "  -1.60%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "65.24"
$ws.Cells.Item(45, 5).Value = "  -1.34%  "
$ws.Cells.Item(46, 5).Value = "  +1.49%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.706"
$ws.Cells.Item(47, 5).Value = "  -1.00%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "6.934"
$ws.Cells.Item(48, 5).Value = "  -1.98%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.1139"
$ws.Cells.Item(49, 5).Value = "  -3.20%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.860"
$ws.Cells.Item(50, 5).Value = "  -1.94%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.3899"
$ws.Cells.Item(51, 5).Value = "  -1.02%  "
